# Applies the "cryptos" price/volume refresh described by the commit diff.
# Each target cell currently holds an inline/shared STRING (not a number) even
# though many values look numeric (e.g. "0.999", "65.948.30"). Setting .Value
# directly on such a string lets Excel auto-convert it to a real number, which
# would change the stored cell type. To keep the cell type as text (matching the
# original workbook), we briefly force a Text number format ("@") before writing
# the value, then restore the cell style to "Normal" so no stray formatting is
# left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "65.948.30"
Set-TextValue "E2" "  +0.73%  "

# Row 3
Set-TextValue "D3" "2.666.04"
Set-TextValue "E3" "  +0.19%  "

# Row 4
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.11%  "

# Row 5
Set-TextValue "D5" "600.77"
Set-TextValue "E5" "  +0.60%  "

# Row 6
Set-TextValue "D6" "161.30"
Set-TextValue "E6" "  +3.43%  "

# Row 7
Set-TextValue "E7" "  +4.77%  "

# Row 8
Set-TextValue "D8" "0.999"
Set-TextValue "E8" "  -0.19%  "

# Row 9
Set-TextValue "D9" "0.128"
Set-TextValue "E9" "  -0.12%  "

# Row 10
Set-TextValue "E10" "  +1.30%  "

# Row 11
Set-TextValue "D11" "5.90"
Set-TextValue "E11" "  +1.19%  "

# Row 12
Set-TextValue "E12" "  +1.79%  "

# Row 13
Set-TextValue "D13" "29.34"
Set-TextValue "E13" "  +0.63%  "

# Row 14
Set-TextValue "D14" "0.0000197"
Set-TextValue "E14" "  +1.28%  "

# Row 15
Set-TextValue "D15" "3.144.43"
Set-TextValue "E15" "  +0.05%  "

# Row 16
Set-TextValue "D16" "65.805.20"
Set-TextValue "E16" "  +0.66%  "

# Row 17
Set-TextValue "D17" "2.680.46"
Set-TextValue "E17" "  -0.05%  "

# Row 18
Set-TextValue "D18" "12.67"
Set-TextValue "E18" "  -0.69%  "

# Row 19
Set-TextValue "E19" "  +1.36%  "

# Row 20
Set-TextValue "D20" "358.00"
Set-TextValue "E20" "  +2.30%  "

# Row 21
Set-TextValue "D21" "7.52"
Set-TextValue "E21" "  +0.20%  "

# Row 22
Set-TextValue "E22" "  -0.05%  "

# Row 23
Set-TextValue "D23" "70.13"
Set-TextValue "E23" "  +0.82%  "

# Row 24
Set-TextValue "E24" "  +11.48%  "

# Row 25
Set-TextValue "D25" "0.0000115"
Set-TextValue "E25" "  +4.28%  "

# Row 26
Set-TextValue "D26" "9.82"
Set-TextValue "E26" "  +2.83%  "

# Row 27
Set-TextValue "E27" "  +3.60%  "

# Row 28
Set-TextValue "D28" "580.42"
Set-TextValue "E28" "  +11.57%  "

# Row 29
Set-TextValue "D29" "8.19"
Set-TextValue "E29" "  +2.84%  "

# Row 30
Set-TextValue "E30" "  -0.92%  "

# Row 31
Set-TextValue "D31" "1.00"
Set-TextValue "E31" "  -0.06%  "

# Row 32
Set-TextValue "E32" "  +1.53%  "

# Row 33
Set-TextValue "E33" "  +4.95%  "

# Row 34
Set-TextValue "D34" "6.78"
Set-TextValue "E34" "  +6.03%  "

# Row 35
Set-TextValue "E35" "  +1.91%  "

# Row 36
Set-TextValue "D36" "0.424"
Set-TextValue "E36" "  +0.90%  "

# Row 37
Set-TextValue "D37" "20.68"
Set-TextValue "E37" "  +0.94%  "

# Row 38
Set-TextValue "B38" "Stacks"
Set-TextValue "C38" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D38" "1.98"
Set-TextValue "E38" "  +3.15%  "

# Row 39
Set-TextValue "B39" "FirstDigitalUSD"
Set-TextValue "C39" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D39" "0.999"
Set-TextValue "E39" "  +0.05%  "

# Row 40
Set-TextValue "D40" "154.32"
Set-TextValue "E40" "  -2.07%  "

# Row 41
Set-TextValue "D41" "2.56"
Set-TextValue "E41" "  +12.23%  "

# Row 42
Set-TextValue "D42" "162.74"
Set-TextValue "E42" "  +0.20%  "

# Row 43
Set-TextValue "D43" "4.13"
Set-TextValue "E43" "  +0.96%  "

# Row 44
Set-TextValue "D44" "0.0622"
Set-TextValue "E44" "  +3.10%  "

# Row 45
Set-TextValue "D45" "23.66"
Set-TextValue "E45" "  +4.85%  "

# Row 46
Set-TextValue "D46" "0.647"
Set-TextValue "E46" "  +1.55%  "

# Row 47
Set-TextValue "E47" "  +1.72%  "

# Row 48
Set-TextValue "E48" "  +2.45%  "

# Row 49
Set-TextValue "D49" "19.86"
Set-TextValue "E49" "  -0.31%  "

# Row 50
Set-TextValue "D50" "0.0₆0248"
Set-TextValue "E50" "  -4.91%  "

# Row 51
Set-TextValue "D51" "0.822"
Set-TextValue "E51" "  +2.29%  "
